$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.761402010917664
$ws.Range("B1").Value = 3.880680084228516
$ws.Range("C1").Value = 5.835652351379395
$ws.Range("D1").Value = 1.546212792396545
$ws.Range("E1").Value = 0.8441539406776428
